# Fuel Prod Imp Exp Balancing Priorities.xlsx - "CA 4.0 files test" edit
#
# This script reproduces the data/formatting changes from the commit:
#  - About!C1: refresh the "last updated" date
#  - About!A10, A11, A37: drop a redundant (no-op) font override style
#  - FPIEBP: drop redundant (no-op) font-override styles on column A labels
#    and on several B:D priority cells (and remove the now-unused, always-
#    empty E9 cell, which shrinks the sheet's used range to A1:D22)
#  - FPIEBP: swap the priority-1/priority-2 ranking for six fuels
#    (natural gas, petroleum gasoline, petroleum diesel, jet fuel/kerosene,
#    heavy fuel oil, LPG propane or butane)
#  - FPIEBP: update the saved selection/active cell

$wb = $excel.ActiveWorkbook

$about = $wb.Worksheets.Item("About")
$fpiebp = $wb.Worksheets.Item("FPIEBP")

# ---------------------------------------------------------------------------
# About sheet
# ---------------------------------------------------------------------------

# Updated "last revised" date (serial 44536 -> 44872, i.e. 12/6/2021 -> 11/7/2022)
$about.Range("C1").Value = 44872

# These three labels had a vestigial "applyFont" style (font id 0 == default
# font, so visually identical to no style at all); clear it so the cell goes
# back to the workbook's default/unstyled state.
$about.Range("A10").Font.Bold = $false
$about.Range("A11").Font.Bold = $false
$about.Range("A37").Font.Bold = $false

# ---------------------------------------------------------------------------
# FPIEBP sheet
# ---------------------------------------------------------------------------

# Column-A row labels: same vestigial no-op font style cleared as above.
foreach ($row in 3,4,5,9,10,11,12,13,14,17,18,19,20,21,22) {
    $fpiebp.Cells.Item($row, 1).Font.Bold = $false
}

# B:D priority cells on the "real fuel" rows had a vestigial no-op fill
# style (fill id 0 == default/no fill); clear it the same way.
foreach ($row in 9,10,11,12,13,14,17,18,19,20,21,22) {
    foreach ($col in 2,3,4) {
        $fpiebp.Cells.Item($row, $col).Font.Bold = $false
    }
}

# The always-empty E9 cell is removed entirely (no longer part of the used
# range); this also shrinks the sheet dimension from A1:E22 to A1:D22.
$fpiebp.Range("E9").Clear()

# Priority swaps: move each of these six fuels from priority order
# (production=3, imports=1, exports=2) to (production=1, imports=3, exports=2).
# natural gas
$fpiebp.Range("B4").Value = 3
$fpiebp.Range("C4").Value = 1
# petroleum gasoline
$fpiebp.Range("B10").Value = 1
$fpiebp.Range("C10").Value = 3
# petroleum diesel
$fpiebp.Range("B11").Value = 1
$fpiebp.Range("C11").Value = 3
# jet fuel or kerosene
$fpiebp.Range("B14").Value = 1
$fpiebp.Range("C14").Value = 3
# heavy fuel oil
$fpiebp.Range("B19").Value = 1
$fpiebp.Range("C19").Value = 3
# LPG propane or butane
$fpiebp.Range("B20").Value = 1
$fpiebp.Range("C20").Value = 3

# Update the saved selection to match (FPIEBP tab active, G6 selected).
$fpiebp.Activate()
$fpiebp.Range("G6").Select()
